$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Förändrad" date column (C2:C33) from 45183 to 45184 ---
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# --- Update the HYPERLINK formulas in row 2 (columns S, T, V, W, X, Y) to include
#     the friendly-name argument "A 33036-2023" ---

# S2 keeps the (malformed, as supplied) target formula text exactly as given.
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/artfynd/A 33036-2023.xlsx, "A 33036-2023"")'

$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/kartor/A 33036-2023.png", "A 33036-2023")'

$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/klagomål/A 33036-2023.docx", "A 33036-2023")'

$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/klagomålsmail/A 33036-2023.docx", "A 33036-2023")'

$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/tillsyn/A 33036-2023.docx", "A 33036-2023")'

# Y2 was stored as a plain inline-string (not a real formula) before the edit;
# it becomes an actual formula now.
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/tillsynsmail/A 33036-2023.docx", "A 33036-2023")'

Write-Output "edit applied"
